$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(5, 7).Value = 459596.8893580333
$ws.Cells.Item(5, 8).Value = 87440.32838484559
$ws.Cells.Item(5, 10).Value = 32376.06915003818

$ws.Cells.Item(6, 7).Value = 459596.8893580333
$ws.Cells.Item(6, 8).Value = 87440.32838484559
$ws.Cells.Item(6, 10).Value = 32376.06915003818

$ws.Cells.Item(7, 7).Value = 919193.7787160666
$ws.Cells.Item(7, 8).Value = 175939.9156559471
$ws.Cells.Item(7, 10).Value = 65281.23913692087

$ws.Cells.Item(8, 7).Value = 919193.7787160666
$ws.Cells.Item(8, 8).Value = 175939.9156559471
$ws.Cells.Item(8, 10).Value = 65281.23913692087

$ws.Cells.Item(9, 7).Value = 1378790.6680741
$ws.Cells.Item(9, 8).Value = 265499.7793553168
$ws.Cells.Item(9, 10).Value = 98567.32161440127

$ws.Cells.Item(10, 7).Value = 1378790.6680741
$ws.Cells.Item(10, 8).Value = 265499.7793553168
$ws.Cells.Item(10, 10).Value = 98567.32161440127

$ws.Cells.Item(11, 7).Value = 1838387.557432133
$ws.Cells.Item(11, 8).Value = 351604.8140642748
$ws.Cells.Item(11, 9).Value = 0.873350509689706
$ws.Cells.Item(11, 10).Value = 130857.7343420096

$ws.Cells.Item(12, 7).Value = 1838387.557432133
$ws.Cells.Item(12, 8).Value = 351604.8140642748
$ws.Cells.Item(12, 9).Value = 0.873350509689706
$ws.Cells.Item(12, 10).Value = 130857.7343420096

$ws.Cells.Item(13, 7).Value = 2297984.446790164
$ws.Cells.Item(13, 8).Value = 444388.4130639153
$ws.Cells.Item(13, 10).Value = 165013.6639672196

$ws.Cells.Item(14, 7).Value = 2297984.446790164
$ws.Cells.Item(14, 8).Value = 444388.4130639153
$ws.Cells.Item(14, 10).Value = 165013.6639672196

$ws.Cells.Item(15, 7).Value = 2757581.3361482
$ws.Cells.Item(15, 8).Value = 527519.8281481991
$ws.Cells.Item(15, 9).Value = 1.310025764534559
$ws.Cells.Item(15, 10).Value = 195382.086049675

$ws.Cells.Item(16, 7).Value = 2757581.3361482
$ws.Cells.Item(16, 8).Value = 527519.8281481991
$ws.Cells.Item(16, 9).Value = 1.310025764534559
$ws.Cells.Item(16, 10).Value = 195382.086049675

$ws.Cells.Item(17, 7).Value = 2757581.3361482
$ws.Cells.Item(17, 8).Value = 527519.8281481991
$ws.Cells.Item(17, 9).Value = 1.310025764534559
$ws.Cells.Item(17, 10).Value = 195382.086049675

$ws.Cells.Item(18, 7).Value = 2757581.3361482
$ws.Cells.Item(18, 8).Value = 527519.8281481991
$ws.Cells.Item(18, 9).Value = 1.310025764534559
$ws.Cells.Item(18, 10).Value = 195382.086049675

$ws.Cells.Item(19, 7).Value = 2757581.3361482
$ws.Cells.Item(19, 8).Value = 527519.8281481991
$ws.Cells.Item(19, 9).Value = 1.310025764534559
$ws.Cells.Item(19, 10).Value = 195382.086049675

$ws.Cells.Item(20, 7).Value = 2757581.3361482
$ws.Cells.Item(20, 8).Value = 527519.8281481991
$ws.Cells.Item(20, 9).Value = 1.310025764534559
$ws.Cells.Item(20, 10).Value = 195382.086049675

$ws.Cells.Item(21, 7).Value = 2757581.3361482
$ws.Cells.Item(21, 8).Value = 527519.8281481991
$ws.Cells.Item(21, 9).Value = 1.310025764534559
$ws.Cells.Item(21, 10).Value = 195382.086049675

$ws.Cells.Item(22, 7).Value = 2757581.3361482
$ws.Cells.Item(22, 8).Value = 527519.8281481991
$ws.Cells.Item(22, 9).Value = 1.310025764534559
$ws.Cells.Item(22, 10).Value = 195382.086049675

$ws.Cells.Item(23, 7).Value = 2757581.3361482
$ws.Cells.Item(23, 8).Value = 527519.8281481991
$ws.Cells.Item(23, 9).Value = 1.310025764534559
$ws.Cells.Item(23, 10).Value = 195382.086049675

$ws.Cells.Item(24, 7).Value = 2757581.3361482
$ws.Cells.Item(24, 8).Value = 527519.8281481991
$ws.Cells.Item(24, 9).Value = 1.310025764534559
$ws.Cells.Item(24, 10).Value = 195382.086049675

$ws.Cells.Item(25, 7).Value = 2757581.3361482
$ws.Cells.Item(25, 8).Value = 527519.8281481991
$ws.Cells.Item(25, 9).Value = 1.310025764534559
$ws.Cells.Item(25, 10).Value = 195382.086049675

$ws.Cells.Item(26, 7).Value = 2757581.3361482
$ws.Cells.Item(26, 8).Value = 527519.8281481991
$ws.Cells.Item(26, 9).Value = 1.310025764534559
$ws.Cells.Item(26, 10).Value = 195382.086049675

$ws.Cells.Item(27, 7).Value = 2757581.3361482
$ws.Cells.Item(27, 8).Value = 527519.8281481991
$ws.Cells.Item(27, 9).Value = 1.310025764534559
$ws.Cells.Item(27, 10).Value = 195382.086049675

$ws.Cells.Item(28, 7).Value = 2757581.3361482
$ws.Cells.Item(28, 8).Value = 527519.8281481991
$ws.Cells.Item(28, 9).Value = 1.310025764534559
$ws.Cells.Item(28, 10).Value = 195382.086049675

$ws.Cells.Item(29, 7).Value = 2757581.3361482
$ws.Cells.Item(29, 8).Value = 527519.8281481991
$ws.Cells.Item(29, 9).Value = 1.310025764534559
$ws.Cells.Item(29, 10).Value = 195382.086049675

$ws.Cells.Item(30, 7).Value = 2757581.3361482
$ws.Cells.Item(30, 8).Value = 527519.8281481991
$ws.Cells.Item(30, 9).Value = 1.310025764534559
$ws.Cells.Item(30, 10).Value = 195382.086049675

$ws.Cells.Item(31, 7).Value = 2757581.3361482
$ws.Cells.Item(31, 8).Value = 527519.8281481991
$ws.Cells.Item(31, 9).Value = 1.310025764534559
$ws.Cells.Item(31, 10).Value = 195382.086049675

$ws.Cells.Item(32, 7).Value = 2757581.3361482
$ws.Cells.Item(32, 8).Value = 527519.8281481991
$ws.Cells.Item(32, 9).Value = 1.310025764534559
$ws.Cells.Item(32, 10).Value = 195382.086049675

$ws.Cells.Item(33, 7).Value = 2757581.3361482
$ws.Cells.Item(33, 8).Value = 527519.8281481991
$ws.Cells.Item(33, 9).Value = 1.310025764534559
$ws.Cells.Item(33, 10).Value = 195382.086049675

$ws.Cells.Item(34, 7).Value = 2757581.3361482
$ws.Cells.Item(34, 8).Value = 527519.8281481991
$ws.Cells.Item(34, 9).Value = 1.310025764534559
$ws.Cells.Item(34, 10).Value = 195382.086049675

$ws.Cells.Item(35, 7).Value = 2757581.3361482
$ws.Cells.Item(35, 8).Value = 527519.8281481991
$ws.Cells.Item(35, 9).Value = 1.310025764534559
$ws.Cells.Item(35, 10).Value = 195382.086049675

$ws.Cells.Item(36, 7).Value = 2757581.3361482
$ws.Cells.Item(36, 8).Value = 527519.8281481991
$ws.Cells.Item(36, 9).Value = 1.310025764534559
$ws.Cells.Item(36, 10).Value = 195382.086049675

$ws.Cells.Item(37, 7).Value = 2757581.3361482
$ws.Cells.Item(37, 8).Value = 527519.8281481991
$ws.Cells.Item(37, 9).Value = 1.310025764534559
$ws.Cells.Item(37, 10).Value = 195382.086049675

$ws.Cells.Item(38, 7).Value = 2757581.3361482
$ws.Cells.Item(38, 8).Value = 527519.8281481991
$ws.Cells.Item(38, 9).Value = 1.310025764534559
$ws.Cells.Item(38, 10).Value = 195382.086049675

$ws.Cells.Item(39, 7).Value = 2757581.3361482
$ws.Cells.Item(39, 8).Value = 527519.8281481991
$ws.Cells.Item(39, 9).Value = 1.310025764534559
$ws.Cells.Item(39, 10).Value = 195382.086049675

$ws.Cells.Item(40, 7).Value = 2757581.3361482
$ws.Cells.Item(40, 8).Value = 527519.8281481991
$ws.Cells.Item(40, 9).Value = 1.310025764534559
$ws.Cells.Item(40, 10).Value = 195382.086049675

$ws.Cells.Item(41, 7).Value = 2757581.3361482
$ws.Cells.Item(41, 8).Value = 527519.8281481991
$ws.Cells.Item(41, 9).Value = 1.310025764534559
$ws.Cells.Item(41, 10).Value = 195382.086049675

$ws.Cells.Item(42, 7).Value = 2757581.3361482
$ws.Cells.Item(42, 8).Value = 527519.8281481991
$ws.Cells.Item(42, 9).Value = 1.310025764534559
$ws.Cells.Item(42, 10).Value = 195382.086049675

$ws.Cells.Item(43, 7).Value = 2757581.3361482
$ws.Cells.Item(43, 8).Value = 527519.8281481991
$ws.Cells.Item(43, 9).Value = 1.310025764534559
$ws.Cells.Item(43, 10).Value = 195382.086049675

$ws.Cells.Item(44, 7).Value = 2757581.3361482
$ws.Cells.Item(44, 8).Value = 527519.8281481991
$ws.Cells.Item(44, 9).Value = 1.310025764534559
$ws.Cells.Item(44, 10).Value = 195382.086049675

$ws.Cells.Item(45, 7).Value = 2757581.3361482
$ws.Cells.Item(45, 8).Value = 527519.8281481991
$ws.Cells.Item(45, 9).Value = 1.310025764534559
$ws.Cells.Item(45, 10).Value = 195382.086049675

$ws.Cells.Item(46, 7).Value = 2757581.3361482
$ws.Cells.Item(46, 8).Value = 527519.8281481991
$ws.Cells.Item(46, 9).Value = 1.310025764534559
$ws.Cells.Item(46, 10).Value = 195382.086049675

$ws.Cells.Item(47, 7).Value = 2757581.3361482
$ws.Cells.Item(47, 8).Value = 527519.8281481991
$ws.Cells.Item(47, 9).Value = 1.310025764534559
$ws.Cells.Item(47, 10).Value = 195382.086049675

$ws.Cells.Item(48, 7).Value = 2757581.3361482
$ws.Cells.Item(48, 8).Value = 527519.8281481991
$ws.Cells.Item(48, 9).Value = 1.310025764534559
$ws.Cells.Item(48, 10).Value = 195382.086049675

$ws.Cells.Item(49, 7).Value = 2757581.3361482
$ws.Cells.Item(49, 8).Value = 527519.8281481991
$ws.Cells.Item(49, 9).Value = 1.310025764534559
$ws.Cells.Item(49, 10).Value = 195382.086049675

$ws.Cells.Item(50, 7).Value = 2757581.3361482
$ws.Cells.Item(50, 8).Value = 527519.8281481991
$ws.Cells.Item(50, 9).Value = 1.310025764534559
$ws.Cells.Item(50, 10).Value = 195382.086049675

$ws.Cells.Item(51, 7).Value = 2757581.3361482
$ws.Cells.Item(51, 8).Value = 527519.8281481991
$ws.Cells.Item(51, 9).Value = 1.310025764534559
$ws.Cells.Item(51, 10).Value = 195382.086049675

$ws.Cells.Item(52, 7).Value = 2757581.3361482
$ws.Cells.Item(52, 8).Value = 527519.8281481991
$ws.Cells.Item(52, 9).Value = 1.310025764534559
$ws.Cells.Item(52, 10).Value = 195382.086049675

$ws.Cells.Item(53, 7).Value = 2757581.3361482
$ws.Cells.Item(53, 8).Value = 527519.8281481991
$ws.Cells.Item(53, 9).Value = 1.310025764534559
$ws.Cells.Item(53, 10).Value = 195382.086049675

$ws.Cells.Item(54, 7).Value = 2757581.3361482
$ws.Cells.Item(54, 8).Value = 527519.8281481991
$ws.Cells.Item(54, 9).Value = 1.310025764534559
$ws.Cells.Item(54, 10).Value = 195382.086049675

$ws.Cells.Item(55, 7).Value = 2757581.3361482
$ws.Cells.Item(55, 8).Value = 527519.8281481991
$ws.Cells.Item(55, 9).Value = 1.310025764534559
$ws.Cells.Item(55, 10).Value = 195382.086049675

$ws.Cells.Item(56, 7).Value = 2757581.3361482
$ws.Cells.Item(56, 8).Value = 527519.8281481991
$ws.Cells.Item(56, 9).Value = 1.310025764534559
$ws.Cells.Item(56, 10).Value = 195382.086049675

$ws.Cells.Item(57, 7).Value = 2757581.3361482
$ws.Cells.Item(57, 8).Value = 527519.8281481991
$ws.Cells.Item(57, 9).Value = 1.310025764534559
$ws.Cells.Item(57, 10).Value = 195382.086049675

$ws.Cells.Item(58, 7).Value = 2757581.3361482
$ws.Cells.Item(58, 8).Value = 527519.8281481991
$ws.Cells.Item(58, 9).Value = 1.310025764534559
$ws.Cells.Item(58, 10).Value = 195382.086049675

$ws.Cells.Item(59, 7).Value = 2757581.3361482
$ws.Cells.Item(59, 8).Value = 527519.8281481991
$ws.Cells.Item(59, 9).Value = 1.310025764534559
$ws.Cells.Item(59, 10).Value = 195382.086049675

$ws.Cells.Item(60, 7).Value = 2757581.3361482
$ws.Cells.Item(60, 8).Value = 527519.8281481991
$ws.Cells.Item(60, 9).Value = 1.310025764534559
$ws.Cells.Item(60, 10).Value = 195382.086049675

$ws.Cells.Item(61, 7).Value = 2757581.3361482
$ws.Cells.Item(61, 8).Value = 527519.8281481991
$ws.Cells.Item(61, 9).Value = 1.310025764534559
$ws.Cells.Item(61, 10).Value = 195382.086049675

$ws.Cells.Item(62, 7).Value = 2757581.3361482
$ws.Cells.Item(62, 8).Value = 527519.8281481991
$ws.Cells.Item(62, 9).Value = 1.310025764534559
$ws.Cells.Item(62, 10).Value = 195382.086049675

$ws.Cells.Item(63, 7).Value = 2757581.3361482
$ws.Cells.Item(63, 8).Value = 527519.8281481991
$ws.Cells.Item(63, 9).Value = 1.310025764534559
$ws.Cells.Item(63, 10).Value = 195382.086049675

$ws.Cells.Item(64, 7).Value = 2757581.3361482
$ws.Cells.Item(64, 8).Value = 527519.8281481991
$ws.Cells.Item(64, 9).Value = 1.310025764534559
$ws.Cells.Item(64, 10).Value = 195382.086049675

$ws.Cells.Item(65, 7).Value = 2757581.3361482
$ws.Cells.Item(65, 8).Value = 527519.8281481991
$ws.Cells.Item(65, 9).Value = 1.310025764534559
$ws.Cells.Item(65, 10).Value = 195382.086049675

$ws.Cells.Item(66, 7).Value = 2757581.3361482
$ws.Cells.Item(66, 8).Value = 527519.8281481991
$ws.Cells.Item(66, 9).Value = 1.310025764534559
$ws.Cells.Item(66, 10).Value = 195382.086049675

$ws.Cells.Item(67, 7).Value = 2757581.3361482
$ws.Cells.Item(67, 8).Value = 527519.8281481991
$ws.Cells.Item(67, 9).Value = 1.310025764534559
$ws.Cells.Item(67, 10).Value = 195382.086049675
